$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / mean calculation
$ws.Range("F2").Value = -4
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -2
